$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the three new columns: Wins, Losses, Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) by copying
# the format from the last existing header cell (AC1) onto the new headers.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins/Losses/Ties) for every player data row (2-37)
$ws.Range("AD2:AD37").Value = 86
$ws.Range("AE2:AE37").Value = 76
$ws.Range("AF2:AF37").Value = 0
